# "Colocando header nos graficos" - add descriptive header labels to column A
# (row 1) on each data sheet, drop the old header-style formatting from the
# row-label cells below it, fix a handful of missing Portuguese accents, and
# remove the now-unused "Teto" row from the Emissoes sheet. Also update the
# Custo Total sheet to have a proper per-year layout.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Helper: give a cell the same style as a neighbouring "header" cell
# (copy formats only, then restore the text) so we reuse the existing
# style index (s="1") instead of minting a brand-new one.
# ---------------------------------------------------------------------
function Set-HeaderCell($ws, $addr, $styleSourceAddr, $text) {
    $ws.Range($styleSourceAddr).Copy()
    $ws.Range($addr).PasteSpecial(-4122)
    $ws.Range($addr).Value = $text
}

# ---------------------------------------------------------------------
# Sheets 1-3: "Potencia Acumulada - SIN (MW)", "Geracao Periodo Medio
# (MWMed)", "Atendimento a Ponta(MW)" all share the same row/column
# layout (years across B1:E1, technologies down A2:A12).
# ---------------------------------------------------------------------
$sheetIndexes = @(1, 2, 3)
foreach ($idx in $sheetIndexes) {
    $ws = $wb.Worksheets.Item($idx)

    # New column header in A1, using B1's existing header style.
    Set-HeaderCell $ws "A1" "B1" "Fonte/Tecnologia"

    # Row labels A2:A12 lose the bold/bordered header style (back to
    # Normal/default) and a few get accented correctly.
    $ws.Range("A2").Style = "Normal"
    $ws.Range("A2").Value = "Hidro"

    $ws.Range("A3").Style = "Normal"
    $ws.Range("A3").Value = "Gás Natural"

    $ws.Range("A4").Style = "Normal"
    $ws.Range("A4").Value = "Carvão"

    $ws.Range("A5").Style = "Normal"
    $ws.Range("A5").Value = "Nuclear"

    $ws.Range("A6").Style = "Normal"
    $ws.Range("A6").Value = "Óleos Comb"

    $ws.Range("A7").Style = "Normal"
    $ws.Range("A7").Value = "Biomassa"

    $ws.Range("A8").Style = "Normal"
    $ws.Range("A8").Value = "Eólica"

    $ws.Range("A9").Style = "Normal"
    $ws.Range("A9").Value = "Solar"

    $ws.Range("A10").Style = "Normal"
    $ws.Range("A10").Value = "Outros"

    $ws.Range("A11").Style = "Normal"
    $ws.Range("A11").Value = "Pot. Compl."

    $ws.Range("A12").Style = "Normal"
    $ws.Range("A12").Value = "GD"
}

# ---------------------------------------------------------------------
# Sheet 4 "Potencia Incremental - SIN(MW)" is untouched by this commit.
# ---------------------------------------------------------------------

# ---------------------------------------------------------------------
# Sheet 5 "Emissoes Totais (MtCO2eq)": add "Periodo" header, drop the
# header style + fix accents on the two remaining row labels, and
# delete the obsolete "Teto" row entirely.
# ---------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item(5)

Set-HeaderCell $ws5 "A1" "B1" "Período"

$ws5.Range("A2").Style = "Normal"
$ws5.Range("A2").Value = "P.Médio"

$ws5.Range("A3").Style = "Normal"
$ws5.Range("A3").Value = "P.Crítico"

$ws5.Rows.Item(4).Delete()

# ---------------------------------------------------------------------
# Sheet 6 "Custo Total (bilhoes de R$)": turn the single "Custo" column
# into a proper "2015" data column with a new row-label header, fix
# accents, and update the two cost figures.
# ---------------------------------------------------------------------
$ws6 = $wb.Worksheets.Item(6)

Set-HeaderCell $ws6 "A1" "B1" "Tipo Expansão"

# "2015" must stay a text label (like every other sheet's year header),
# not become a number -- paste the value only from an existing "2015"
# text cell so it keeps B1's current style (s="1") untouched.
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("B1").Copy()
$ws6.Range("B1").PasteSpecial(-4163)

$ws6.Range("A2").Style = "Normal"
$ws6.Range("A2").Value = "Expansão Centralizada"
$ws6.Range("B2").Value = 594

$ws6.Range("A3").Style = "Normal"
$ws6.Range("A3").Value = "Expansão por GD"
$ws6.Range("B3").Value = 99
